$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Jaden Ivey"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Detroit Pistons"

# Row 7
$ws.Range("A7").Value = "Giannis Antetokounmpo"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Milwaukee Bucks"

# Row 8
$ws.Range("A8").Value = "Andrew Wiggins"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Golden State Warriors"

# Row 9
$ws.Range("A9").Value = "Jaren Jackson Jr."
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Memphis Grizzlies"

# Row 13
$ws.Range("A13").Value = "Anfernee Simons"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Portland Trail Blazers"

# Row 14
$ws.Range("A14").Value = "Payton Pritchard"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Boston Celtics"

# Row 15
$ws.Range("A15").Value = "Kyle Kuzma"
$ws.Range("B15").Value = "PF"
$ws.Range("C15").Value = "Washington Wizards"

# Row 16
$ws.Range("A16").Value = "Anthony Edwards"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Minnesota Timberwolves"

# Row 17
$ws.Range("A17").Value = "Wendell Carter Jr."
$ws.Range("B17").Value = "C"
$ws.Range("C17").Value = "Orlando Magic"

# Row 18
$ws.Range("A18").Value = "Aaron Gordon"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Denver Nuggets"
